$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the entire "is_viewed" column (column I)
$ws.Columns.Item(9).Delete()

# Update subcategory labels to pluralized/parenthesized form
$rowsToUpdate = @(3, 5, 6, 7, 8, 9, 20, 22)
foreach ($r in $rowsToUpdate) {
    $ws.Cells.Item($r, 8).Value = "line graph(s)"
}

$ws.Cells.Item(23, 8).Value = "bar chart(s)"
$ws.Cells.Item(32, 8).Value = "line graph(s)"
